$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.851.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.373.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.17%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.598"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.363.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.189"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.595"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.900.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "606.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -10.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.630.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.74%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.356.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.84%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.118"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.913"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.28%  "
$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.31%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.63%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.872.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.106"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "529.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.55%  "
$ws.Range("B40").Value = "CoreDAO"
$ws.Range("C40").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +38.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0727"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.128"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.352"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0421"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.131"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.44%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.62%  "
